$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the two "Description" values that are now prefixed with "to_"
#    (Customers -> to_Customers, Services -> to_Services). Doing this FIRST
#    (and before anything else references these strings) controls where the
#    renamed strings land once the shared-string table is recompacted on save.
# ---------------------------------------------------------------------------
$customerRows = @(6, 11, 16, 21, 26, 31, 36, 41, 46)
foreach ($r in $customerRows) {
    $ws.Cells.Item($r, 7).Value = "to_Customers"
}
$ws.Cells.Item(77, 7).Value = "to_Services"

# ---------------------------------------------------------------------------
# 2) Add the new "Type" column (column I) with its header.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 9).Value = "Type"

# ---------------------------------------------------------------------------
# 3) Fill in column I ("Type") for every data row. The per-row classification:
#      lo0        -> loopback
#      ge-0/0/1   -> access     (on MA routers)
#      ge-0/0/0   -> access     (on MA routers) / aggregation (on PAG/AG)
#      ge-0/0/2   -> customer
#      other PAG/AG links -> aggregation / core
#      Core-RR service link -> service
#    Write the NON-loopback rows first (in row order) and the loopback rows
#    last, so that the new unique strings are appended to the shared-string
#    table in the order: access, customer, aggregation, core, service, loopback.
# ---------------------------------------------------------------------------
$typeMap = [ordered]@{
    3  = "loopback"
    4  = "access"
    5  = "access"
    6  = "customer"
    8  = "loopback"
    9  = "access"
    10 = "access"
    11 = "customer"
    13 = "loopback"
    14 = "access"
    15 = "access"
    16 = "customer"
    18 = "loopback"
    19 = "access"
    20 = "access"
    21 = "customer"
    23 = "loopback"
    24 = "access"
    25 = "access"
    26 = "customer"
    28 = "loopback"
    29 = "access"
    30 = "access"
    31 = "customer"
    33 = "loopback"
    34 = "access"
    35 = "access"
    36 = "customer"
    38 = "loopback"
    39 = "access"
    40 = "access"
    41 = "customer"
    43 = "loopback"
    44 = "access"
    45 = "access"
    46 = "customer"
    48 = "loopback"
    49 = "access"
    50 = "access"
    51 = "access"
    52 = "aggregation"
    53 = "aggregation"
    55 = "loopback"
    56 = "access"
    57 = "access"
    58 = "access"
    59 = "aggregation"
    60 = "aggregation"
    62 = "loopback"
    63 = "aggregation"
    64 = "aggregation"
    65 = "core"
    66 = "core"
    68 = "loopback"
    69 = "aggregation"
    70 = "aggregation"
    71 = "core"
    72 = "core"
    74 = "loopback"
    75 = "core"
    76 = "core"
    77 = "service"
}

foreach ($r in $typeMap.Keys) {
    if ($typeMap[$r] -ne "loopback") {
        $ws.Cells.Item($r, 9).Value = $typeMap[$r]
    }
}
foreach ($r in $typeMap.Keys) {
    if ($typeMap[$r] -eq "loopback") {
        $ws.Cells.Item($r, 9).Value = $typeMap[$r]
    }
}

# ---------------------------------------------------------------------------
# 4) Move the active selection to I74 (matches the author's final cursor
#    position after filling in the new column).
# ---------------------------------------------------------------------------
$ws.Range("I74").Select()
